# Apply cell updates from the crypto price refresh (GitHub Actions data pull).
# Numeric-looking "Price" values must stay TEXT (the sheet stores these as
# literal strings, e.g. "302.38" or double-dotted "23.058.95"), so cells whose
# new value would otherwise auto-convert to a Number are forced to Text via
# NumberFormat, then the style is put back to Normal so no stray formatting
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.058.95"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.597.03"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3777"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.238"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08111"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.531"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.289"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001234"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "1.595.81"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06845"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.486"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Value = "23.050.47"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.375"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.800"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.228"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.378"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.713"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.19%  "
$ws.Range("D33").Value = "1.771.20"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9566"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07470"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02692"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.129"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08808"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2502"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.362"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7005"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.36%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.005"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.276"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07918"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.229"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.82%  "
